# Update the "Förändrad" (Changed) date column (C) for rows 2-9
# from 45184 (2023-09-15) to 45185 (2023-09-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2:C9").Value = 45185
